# Update workbook per commit: gh-pages data refresh
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("D2").Value = "秋月路9号 南宁五象山庄"
$wsExhibit.Range("F4").Value = 228
$wsExhibit.Range("F5").Value = 2738
$wsExhibit.Range("F6").Value = 1917
$wsExhibit.Range("F8").Value = 120
$wsExhibit.Range("F9").Value = 969
$wsExhibit.Range("F10").Value = 184
$wsExhibit.Range("F11").Value = 13

# Sheet "全部类型" (all types) - has an extra row (concert) inserted at row 8,
# shifting the remaining exhibition rows down by one compared to "展览"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("D2").Value = "秋月路9号 南宁五象山庄"
$wsAll.Range("F4").Value = 228
$wsAll.Range("F5").Value = 2738
$wsAll.Range("F6").Value = 1917
$wsAll.Range("F9").Value = 120
$wsAll.Range("F10").Value = 969
$wsAll.Range("F11").Value = 184
$wsAll.Range("F12").Value = 13
